# Add a new "Swiss" test-data worksheet, modelled on the existing
# "Germany" sheet, for the Switzerland market.

$wb = $excel.ActiveWorkbook

# The Germany sheet is the template: same layout/styles as the other
# per-country sheets, just missing the one extra "MX Minerva Bridge Kit"
# row that Belgium/Czech/Germany carry.
$germany = $wb.Worksheets.Item("Germany")

# Copy it to the end of the workbook - this clones formatting, merged
# cells, column widths, etc. exactly like Excel's "Move or Copy... (Create
# a copy)" command.
$germany.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# This template row ("MX Minerva Bridge Kit") doesn't apply to
# Switzerland, so drop it - everything below shifts up automatically.
$swiss.Rows.Item(11).Delete()

# Fill in the market-specific values.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2343/T2641"

# Restore Germany's own selection/active state (it had the user's
# selection before the copy), then activate the new Swiss tab so it's the
# one shown/selected when the workbook is reopened.
$germany.Select()
$germany.Cells.Select()
$swiss.Select()
